$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh (GitHub Actions bot, 2023-02-08 23:15 UTC):
# several coin rows rotate up by one, prices/volumes update, and every
# row's "Hora" (hour) column moves from 22 to 23.
#
# All of these columns are stored as literal text in the workbook (e.g. a
# price of "327.80" keeps its trailing zero, and a percent like "-1.65%" or
# an hour like "23" are not numeric cells). Excel's default General format
# would happily reinterpret "327.80" / "-1.65%" / "23" as numbers, silently
# dropping the formatting we need to preserve - so each cell is forced to
# Text format before the value is written.

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2
Set-TextValue $ws.Range("D2") '327.80'
Set-TextValue $ws.Range("E2") '-1.65%'
Set-TextValue $ws.Range("G2") '23'

# Row 3
Set-TextValue $ws.Range("D3") '44.29'
Set-TextValue $ws.Range("E3") '-1.03%'
Set-TextValue $ws.Range("G3") '23'

# Row 4
Set-TextValue $ws.Range("D4") '5.373'
Set-TextValue $ws.Range("E4") '-3.30%'
Set-TextValue $ws.Range("G4") '23'

# Row 5
Set-TextValue $ws.Range("D5") '0.08366'
Set-TextValue $ws.Range("E5") '0.45%'
Set-TextValue $ws.Range("G5") '23'

# Row 6
Set-TextValue $ws.Range("B6") 'FTXToken'
Set-TextValue $ws.Range("C6") 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range("D6") '1.943'
Set-TextValue $ws.Range("E6") '-5.06%'
Set-TextValue $ws.Range("G6") '23'

# Row 7
Set-TextValue $ws.Range("B7") 'MXToken'
Set-TextValue $ws.Range("C7") 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D7") '0.9737'
Set-TextValue $ws.Range("E7") '-0.64%'
Set-TextValue $ws.Range("G7") '23'

# Row 8
Set-TextValue $ws.Range("B8") 'BTSEToken'
Set-TextValue $ws.Range("C8") 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range("D8") '2.532'
Set-TextValue $ws.Range("E8") '-3.07%'
Set-TextValue $ws.Range("G8") '23'

# Row 9
Set-TextValue $ws.Range("B9") 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws.Range("C9") 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range("D9") '0.1148'
Set-TextValue $ws.Range("E9") '2.36%'
Set-TextValue $ws.Range("G9") '23'

# Row 10
Set-TextValue $ws.Range("B10") 'WazirX'
Set-TextValue $ws.Range("C10") 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range("D10") '0.1909'
Set-TextValue $ws.Range("E10") '-0.33%'
Set-TextValue $ws.Range("G10") '23'

# Row 11
Set-TextValue $ws.Range("B11") 'MandalaExchangeToken'
Set-TextValue $ws.Range("C11") 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range("D11") '0.09718'
Set-TextValue $ws.Range("E11") '-3.66%'
Set-TextValue $ws.Range("G11") '23'

# Row 12
Set-TextValue $ws.Range("B12") 'BitrueCoin'
Set-TextValue $ws.Range("C12") 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range("D12") '0.04635'
Set-TextValue $ws.Range("E12") '0.10%'
Set-TextValue $ws.Range("G12") '23'

# Row 13
Set-TextValue $ws.Range("B13") 'BitMartToken'
Set-TextValue $ws.Range("C13") 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range("D13") '0.1062'
Set-TextValue $ws.Range("E13") '0.15%'
Set-TextValue $ws.Range("G13") '23'

# Row 14
Set-TextValue $ws.Range("B14") 'BitForexToken'
Set-TextValue $ws.Range("C14") 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range("D14") '0.001293'
Set-TextValue $ws.Range("E14") '1.28%'
Set-TextValue $ws.Range("G14") '23'

# Row 15
Set-TextValue $ws.Range("B15") 'TigerCash'
Set-TextValue $ws.Range("C15") 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range("D15") '0.005852'
Set-TextValue $ws.Range("E15") '-2.39%'
Set-TextValue $ws.Range("G15") '23'

# Row 16
Set-TextValue $ws.Range("B16") 'LEO'
Set-TextValue $ws.Range("C16") 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D16") '3.367'
Set-TextValue $ws.Range("E16") '0.23%'
Set-TextValue $ws.Range("G16") '23'

# Row 17
Set-TextValue $ws.Range("B17") 'GateToken'
Set-TextValue $ws.Range("C17") 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range("D17") '4.438'
Set-TextValue $ws.Range("E17") '0.02%'
Set-TextValue $ws.Range("G17") '23'

# Row 18
Set-TextValue $ws.Range("D18") '0.3357'
Set-TextValue $ws.Range("E18") '0.42%'
Set-TextValue $ws.Range("G18") '23'

# Row 19
Set-TextValue $ws.Range("D19") '8.429'
Set-TextValue $ws.Range("E19") '-18.47%'
Set-TextValue $ws.Range("G19") '23'

# Row 20
Set-TextValue $ws.Range("D20") '0.1352'
Set-TextValue $ws.Range("E20") '-2.43%'
Set-TextValue $ws.Range("G20") '23'

# Row 21
Set-TextValue $ws.Range("E21") '6.57%'
Set-TextValue $ws.Range("G21") '23'

# Row 22
Set-TextValue $ws.Range("D22") '0.04176'
Set-TextValue $ws.Range("E22") '1.31%'
Set-TextValue $ws.Range("G22") '23'

# Row 23
Set-TextValue $ws.Range("D23") '0.001241'
Set-TextValue $ws.Range("E23") '-4.69%'
Set-TextValue $ws.Range("G23") '23'

# Row 24
Set-TextValue $ws.Range("D24") '0.004430'
Set-TextValue $ws.Range("E24") '0.25%'
Set-TextValue $ws.Range("G24") '23'

# Row 25
Set-TextValue $ws.Range("D25") '0.0001302'
Set-TextValue $ws.Range("E25") '1.78%'
Set-TextValue $ws.Range("G25") '23'

# Row 26
Set-TextValue $ws.Range("D26") '0.0002984'
Set-TextValue $ws.Range("E26") '-20.24%'
Set-TextValue $ws.Range("G26") '23'

# Row 27
Set-TextValue $ws.Range("G27") '23'

# Row 28
Set-TextValue $ws.Range("G28") '23'

# Row 29
Set-TextValue $ws.Range("G29") '23'

# Row 30
Set-TextValue $ws.Range("G30") '23'

# Row 31
Set-TextValue $ws.Range("G31") '23'

# Row 32
Set-TextValue $ws.Range("G32") '23'

# Row 33
Set-TextValue $ws.Range("G33") '23'

# Row 34
Set-TextValue $ws.Range("G34") '23'

# Row 35
Set-TextValue $ws.Range("G35") '23'

# Row 36
Set-TextValue $ws.Range("G36") '23'

# Row 37
Set-TextValue $ws.Range("G37") '23'

# Row 38
Set-TextValue $ws.Range("D38") '0.02703'
Set-TextValue $ws.Range("E38") '-4.20%'
Set-TextValue $ws.Range("G38") '23'

# Row 39
Set-TextValue $ws.Range("D39") '0.05620'
Set-TextValue $ws.Range("E39") '-2.79%'
Set-TextValue $ws.Range("G39") '23'

# Row 40
Set-TextValue $ws.Range("D40") '0.007848'
Set-TextValue $ws.Range("E40") '2.59%'
Set-TextValue $ws.Range("G40") '23'

# Row 41
Set-TextValue $ws.Range("D41") '0.1412'
Set-TextValue $ws.Range("E41") '-1.10%'
Set-TextValue $ws.Range("G41") '23'

# Row 42
Set-TextValue $ws.Range("D42") '0.007316'
Set-TextValue $ws.Range("E42") '-3.27%'
Set-TextValue $ws.Range("G42") '23'

# Row 43
Set-TextValue $ws.Range("D43") '0.002054'
Set-TextValue $ws.Range("E43") '4.08%'
Set-TextValue $ws.Range("G43") '23'

# Row 44
Set-TextValue $ws.Range("D44") '0.008721'
Set-TextValue $ws.Range("E44") '8.56%'
Set-TextValue $ws.Range("G44") '23'

# Row 45
Set-TextValue $ws.Range("D45") '0.3507'
Set-TextValue $ws.Range("G45") '23'

# Row 46
Set-TextValue $ws.Range("D46") '0.00006906'
Set-TextValue $ws.Range("E46") '-1.80%'
Set-TextValue $ws.Range("G46") '23'

# Row 47
Set-TextValue $ws.Range("E47") '0.13%'
Set-TextValue $ws.Range("G47") '23'

# Row 48
Set-TextValue $ws.Range("D48") '0.003501'
Set-TextValue $ws.Range("E48") '-0.90%'
Set-TextValue $ws.Range("G48") '23'

# Row 49
Set-TextValue $ws.Range("D49") '0.003536'
Set-TextValue $ws.Range("E49") '40.11%'
Set-TextValue $ws.Range("G49") '23'

# Row 50
Set-TextValue $ws.Range("D50") '0.00002104'
Set-TextValue $ws.Range("E50") '0.13%'
Set-TextValue $ws.Range("G50") '23'

# Row 51
Set-TextValue $ws.Range("D51") '0.0002003'
Set-TextValue $ws.Range("E51") '0.13%'
Set-TextValue $ws.Range("G51") '23'
